# Add "2022-Q4" quarterly data:
#   1. Insert a new worksheet "2022-Q4" before "2022-Q3" (all later sheets just
#      shift right in tab order / rIds, their own content is untouched).
#   2. Populate it with the new quarter's fund holdings.
#   3. Prepend a corresponding summary row on "总计" (shifting the other rows
#      down one, keeping the running 0-based index in column A consistent).
#   4. Restore the originally-selected tab ("2020-Q4") since adding a sheet
#      steals the active-tab flag.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1 & 2: new "2022-Q4" sheet
# ---------------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($sheetQ3)
$newSheet.Name = "2022-Q4"
# Worksheets.Add() seeds the new sheet with a copy of another sheet's data in
# this environment -- wipe it so we start from a blank grid.
$newSheet.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$fundRows = @(
    @("513090", "易方达中证香港证券投资主题ETF",       "11.28", "97.07", "7.35", "0.8291", 5),
    @("003413", "华泰柏瑞新经济沪港深混合",             "1.44",  "94.26", "5.58", "0.0804", 8),
    @("007151", "前海开源沪港深聚瑞混合",               "0.65",  "89.66", "8.12", "0.0528", 2),
    @("011355", "华泰柏瑞港股通时代机遇混合A",           "0.70",  "94.61", "5.50", "0.0385", 9),
    @("011356", "华泰柏瑞港股通时代机遇混合C",           "0.39",  "94.61", "5.50", "0.0214", 9),
    @("001942", "前海开源沪港深汇鑫灵活配置混合A",       "0.31",  "90.30", "4.27", "0.0132", 9),
    @("001943", "前海开源沪港深汇鑫灵活配置混合C",       "0.27",  "90.30", "4.27", "0.0115", 9)
)

$r = 2
foreach ($row in $fundRows) {
    $aCell = $newSheet.Cells.Item($r, 1)
    $aCell.Value = $r - 2
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    # Columns B..G are stored as plain text in the source data (keeps leading
    # zeros on fund codes like "003413" intact), column H is numeric.
    $textRange = $newSheet.Range($newSheet.Cells.Item($r, 2), $newSheet.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    for ($col = 2; $col -le 7; $col++) {
        $newSheet.Cells.Item($r, $col).Value = $row[$col - 2]
    }
    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r++
}

# ---------------------------------------------------------------------------
# 3: prepend summary row on "总计"
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing data rows (2..9) down to (3..10), columns B..D only -- column
# A (the running index) is rewritten from scratch below so it stays 0-based.
for ($row = 9; $row -ge 2; $row--) {
    $dest = $row + 1
    $total.Cells.Item($dest, 2).Value = $total.Cells.Item($row, 2).Value2
    $total.Cells.Item($dest, 3).Value = $total.Cells.Item($row, 3).Value2
    $total.Cells.Item($dest, 4).Value = $total.Cells.Item($row, 4).Value2
}

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 1.05

# Make sure the newly-materialised row 10 carries the same styling as the
# other index cells in column A before the values are rewritten.
$total.Range("A9").Copy()
$total.Range("A10").PasteSpecial(-4122)

for ($row = 2; $row -le 10; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}

# ---------------------------------------------------------------------------
# 4: restore the originally active tab
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
